$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) from "Checklist" to "Session"
$ws.Name = "Session"

# The log row for Student ID 201566 (logged at 12:53:43) was removed from
# the session entirely - find it dynamically and delete the whole row,
# shifting the rows below it up. (Note: use Value2 for reads - Value's
# getter is unreliable in this host, though it works fine for writes.)
$lastRow = $ws.UsedRange.Rows.Count
for ($r = $lastRow; $r -ge 1; $r--) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "201566") {
        $ws.Rows.Item($r).Delete()
    }
}

# Every remaining log entry that was typed as "Selection" is now recorded
# as a "Scan" (column E holds the log Type).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 5).Value2 -eq "Selection") {
        $ws.Cells.Item($r, 5).Value = "Scan"
    }
}

Write-Host "Sheet renamed to '$($ws.Name)'; final used range: $($ws.UsedRange.Address())"
